# Apply the ClimaWatch Defect Report edits:
#  - Fix capitalization of "Javascript" -> "JavaScript" in the Enviroment column (F)
#  - Fix capitalization of "Dbreaver" -> "DBeaver" in the Enviroment column (F)
#  - Update the active selection to F8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defect Report")

$ws.Range("F4").Value = "JavaScript"
$ws.Range("F6").Value = "DBeaver"
$ws.Range("F7").Value = "JavaScript"
$ws.Range("F8").Value = "DBeaver"
$ws.Range("F9").Value = "JavaScript"

$ws.Range("F8").Select()
